# Updated symbol list on Tue Jan 10 05:00:37 UTC 2023 with GitHub Actions
# Refresh Price (D), Volume(1h) (E) and Hora (G) columns for rows 2-51.
# Leading "'" forces the numeric/percent-looking text to stay literal text
# (matching the sheet's existing inlineStr cells) instead of being
# auto-converted to a number/percentage by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.20"
$ws.Range("E2").Value = "'-1.56%"
$ws.Range("G2").Value = "'5"
$ws.Range("D3").Value = "'26.71"
$ws.Range("E3").Value = "'-2.03%"
$ws.Range("G3").Value = "'5"
$ws.Range("D4").Value = "'4.870"
$ws.Range("G4").Value = "'5"
$ws.Range("D5").Value = "'0.06325"
$ws.Range("E5").Value = "'0.92%"
$ws.Range("G5").Value = "'5"
$ws.Range("D6").Value = "'6.884"
$ws.Range("G6").Value = "'5"
$ws.Range("D7").Value = "'3.316"
$ws.Range("E7").Value = "'1.25%"
$ws.Range("G7").Value = "'5"
$ws.Range("D8").Value = "'1.244"
$ws.Range("E8").Value = "'31.89%"
$ws.Range("G8").Value = "'5"
$ws.Range("D9").Value = "'0.8720"
$ws.Range("E9").Value = "'-0.84%"
$ws.Range("G9").Value = "'5"
$ws.Range("D10").Value = "'0.1456"
$ws.Range("E10").Value = "'0.36%"
$ws.Range("G10").Value = "'5"
$ws.Range("D11").Value = "'0.05106"
$ws.Range("E11").Value = "'-0.76%"
$ws.Range("G11").Value = "'5"
$ws.Range("D12").Value = "'0.07367"
$ws.Range("E12").Value = "'1.17%"
$ws.Range("G12").Value = "'5"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("E13").Value = "'-3.54%"
$ws.Range("G13").Value = "'5"
$ws.Range("D14").Value = "'0.09042"
$ws.Range("E14").Value = "'-0.10%"
$ws.Range("G14").Value = "'5"
$ws.Range("E15").Value = "'1.50%"
$ws.Range("G15").Value = "'5"
$ws.Range("D16").Value = "'0.0006290"
$ws.Range("E16").Value = "'0.27%"
$ws.Range("G16").Value = "'5"
$ws.Range("D17").Value = "'0.005964"
$ws.Range("E17").Value = "'-0.23%"
$ws.Range("G17").Value = "'5"
$ws.Range("E18").Value = "'-0.42%"
$ws.Range("G18").Value = "'5"
$ws.Range("E19").Value = "'0.76%"
$ws.Range("G19").Value = "'5"
$ws.Range("G20").Value = "'5"
$ws.Range("D21").Value = "'0.1324"
$ws.Range("E21").Value = "'1.14%"
$ws.Range("G21").Value = "'5"
$ws.Range("D22").Value = "'3.906"
$ws.Range("E22").Value = "'1.62%"
$ws.Range("G22").Value = "'5"
$ws.Range("D23").Value = "'0.04340"
$ws.Range("E23").Value = "'0.38%"
$ws.Range("G23").Value = "'5"
$ws.Range("E24").Value = "'-0.06%"
$ws.Range("G24").Value = "'5"
$ws.Range("D25").Value = "'0.004272"
$ws.Range("E25").Value = "'-0.09%"
$ws.Range("G25").Value = "'5"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("G26").Value = "'5"
$ws.Range("D27").Value = "'0.0001692"
$ws.Range("E27").Value = "'-4.63%"
$ws.Range("G27").Value = "'5"
$ws.Range("G28").Value = "'5"
$ws.Range("G29").Value = "'5"
$ws.Range("G30").Value = "'5"
$ws.Range("G31").Value = "'5"
$ws.Range("G32").Value = "'5"
$ws.Range("G33").Value = "'5"
$ws.Range("G34").Value = "'5"
$ws.Range("G35").Value = "'5"
$ws.Range("G36").Value = "'5"
$ws.Range("G37").Value = "'5"
$ws.Range("G38").Value = "'5"
$ws.Range("G39").Value = "'5"
$ws.Range("E40").Value = "'-0.04%"
$ws.Range("G40").Value = "'5"
$ws.Range("D41").Value = "'0.006726"
$ws.Range("E41").Value = "'0.28%"
$ws.Range("G41").Value = "'5"
$ws.Range("E42").Value = "'1.39%"
$ws.Range("G42").Value = "'5"
$ws.Range("D43").Value = "'0.002098"
$ws.Range("E43").Value = "'-0.14%"
$ws.Range("G43").Value = "'5"
$ws.Range("D44").Value = "'0.01255"
$ws.Range("E44").Value = "'-10.66%"
$ws.Range("G44").Value = "'5"
$ws.Range("D45").Value = "'0.00005320"
$ws.Range("E45").Value = "'2.75%"
$ws.Range("G45").Value = "'5"
$ws.Range("D46").Value = "'2.360"
$ws.Range("E46").Value = "'0.95%"
$ws.Range("G46").Value = "'5"
$ws.Range("D47").Value = "'0.01998"
$ws.Range("E47").Value = "'-33.11%"
$ws.Range("G47").Value = "'5"
$ws.Range("G48").Value = "'5"
$ws.Range("G49").Value = "'5"
$ws.Range("G50").Value = "'5"
$ws.Range("G51").Value = "'5"
